# unitTest_LocalDb.xlsx edit
#
# 1. "#system" sheet, column G ("desktop" named range) gets a new entry
#    "assertElementNotPresent(name)" inserted alphabetically right before the
#    existing "assertElementPresent(name)" row (G5), pushing every row from
#    G5..G97 down by one (new last row becomes G98). Only column G shifts -
#    every other column on those rows is untouched.
# 2. The "desktop" defined name is widened from $G$2:$G$97 to $G$2:$G$98 to
#    cover the newly added row.
# 3. "#system" sheet, cell Y93 (part of the "web" named range) is renamed
#    from "saveISTDivsAsCsv(config,file)" to
#    "saveInfiniteDivsAsCsv(config,file)" - same alphabetical slot, so no
#    shifting needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

$desktopCol = 7   # column G
$webCol = 25      # column Y

# --- shift column G (desktop) down by one starting at row 5, bottom-up ---
for ($r = 97; $r -ge 5; $r--) {
    $v = $ws.Cells.Item($r, $desktopCol).Value2
    $ws.Cells.Item($r + 1, $desktopCol).Value2 = $v
}

# --- insert the new command name into the now-vacated row 5 ---
$ws.Cells.Item(5, $desktopCol).Value2 = "assertElementNotPresent(name)"

# --- widen the "desktop" defined name to include the new row ---
$desktopName = $wb.Names.Item("desktop")
$desktopName.RefersTo = "='#system'!`$G`$2:`$G`$98"

# --- rename saveISTDivsAsCsv(config,file) -> saveInfiniteDivsAsCsv(config,file) ---
$ws.Cells.Item(93, $webCol).Value2 = "saveInfiniteDivsAsCsv(config,file)"
